# Adaption of Excels (remove whitespaces in column GapType!)
#
# The "Gap1_type" column (column M) contains values such as "Arbeit, Privat".
# Remove the stray whitespace after the comma so it reads "Arbeit,Privat".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QAGlist_Teil1")

# Locate the "Gap1_type" column by reading the header row (row 1).
$lastCol = $ws.UsedRange.Columns.Count
$gapTypeCol = 0
for ($c = 1; $c -le $lastCol; $c++) {
    $headerValue = $ws.Cells.Item(1, $c).Value()
    if ($headerValue -eq "Gap1_type") {
        $gapTypeCol = $c
        break
    }
}

if ($gapTypeCol -eq 0) {
    $gapTypeCol = 13 # fallback: column M
}

# Find the last used row of the sheet.
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $gapTypeCol)
    $val = $cell.Value()
    if ($val -ne $null -and $val -is [string] -and $val.Contains(", ")) {
        $cell.Value = $val.Replace(", ", ",")
    }
}

$wb.Save()
